$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32; existing rows 32-53 shift down to 33-54.
$ws.Rows(32).Insert()

# Populate the newly inserted row 32 with its data.
$ws.Cells.Item(32, 1).Value = 10
$ws.Cells.Item(32, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(32, 3).Value = "La Araucanía"
$ws.Cells.Item(32, 4).Value = 44767
$ws.Cells.Item(32, 5).Value = 9
$ws.Cells.Item(32, 6).Value = 100112010
$ws.Cells.Item(32, 7).Value = "Achicoria"
$ws.Cells.Item(32, 8).Value = "Sin especificar"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 500
$ws.Cells.Item(32, 11).Value = 10000
$ws.Cells.Item(32, 12).Value = 11000
$ws.Cells.Item(32, 13).Value = 10600
$ws.Cells.Item(32, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(32, 15).Value = "Región Metropolitana"
$ws.Cells.Item(32, 16).Value = 589
$ws.Cells.Item(32, 17).Value = 18
$ws.Cells.Item(32, 18).Value = "Hortaliza"
